# Apply odds updates for Jogos_da_Semana_FlashScore_2024-10-17
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 1.62
$ws.Range("H2").Value = 3.8
$ws.Range("I2").Value = 5.5
$ws.Range("M2").Value = 1.07
$ws.Range("O2").Value = 1.33
$ws.Range("S2").Value = 1.36
$ws.Range("T2").Value = 3
$ws.Range("W2").Value = 7
$ws.Range("X2").Value = 7.5
$ws.Range("AA2").Value = 13
$ws.Range("AB2").Value = 26
$ws.Range("AC2").Value = 11
$ws.Range("AD2").Value = 7.5
$ws.Range("AE2").Value = 17
$ws.Range("AF2").Value = 51
$ws.Range("AG2").Value = 301
$ws.Range("AM2").Value = 41
$ws.Range("AN2").Value = 3.6
$ws.Range("AP2").Value = 19
$ws.Range("AQ2").Value = 26
$ws.Range("AT2").Value = 3
$ws.Range("AU2").Value = 8.5
$ws.Range("AV2").Value = 51
$ws.Range("AY2").Value = 34
$ws.Range("BA2").Value = 126
$ws.Range("BB2").Value = 251
# Row 3
$ws.Range("M3").Value = 1.11
$ws.Range("N3").Value = 6.5
$ws.Range("P3").Value = 2.5
$ws.Range("Q3").Value = 2.6
$ws.Range("R3").Value = 1.48
$ws.Range("S3").Value = 1.62
$ws.Range("T3").Value = 2.2
$ws.Range("U3").Value = 2.25
$ws.Range("AP3").Value = 26
$ws.Range("AR3").Value = 67
$ws.Range("AS3").Value = 251
$ws.Range("AT3").Value = 2.25
$ws.Range("AU3").Value = 9.5
$ws.Range("AZ3").Value = 101
$ws.Range("BB3").Value = 451
# Row 4
$ws.Range("M4").Value = 1.11
$ws.Range("N4").Value = 6.5
$ws.Range("O4").Value = 1.53
$ws.Range("P4").Value = 2.38
$ws.Range("AQ4").Value = 41
# Row 5
$ws.Range("M5").Value = 1.05
$ws.Range("O5").Value = 1.3
# Row 6
$ws.Range("G6").Value = 1.31
$ws.Range("H6").Value = 4.65
$ws.Range("K6").Value = 2.47
$ws.Range("N6").Value = 13.8
$ws.Range("O6").Value = 1.13
$ws.Range("P6").Value = 4.4
$ws.Range("Q6").Value = 1.52
$ws.Range("R6").Value = 2.2
$ws.Range("S6").Value = 1.27
$ws.Range("T6").Value = 3.42
$ws.Range("U6").Value = 1.84
$ws.Range("V6").Value = 1.92
$ws.Range("W6").Value = 6.8
$ws.Range("X6").Value = 6
$ws.Range("Z6").Value = 7.3
$ws.Range("AA6").Value = 8.75
$ws.Range("AB6").Value = 19.5
$ws.Range("AC6").Value = 14
$ws.Range("AD6").Value = 8.25
$ws.Range("AF6").Value = 60
$ws.Range("AO6").Value = 5.7
$ws.Range("AP6").Value = 14.5
$ws.Range("AQ6").Value = 14.5
$ws.Range("AR6").Value = 37
$ws.Range("AS6").Value = 175
$ws.Range("AT6").Value = 3.3
$ws.Range("AW6").Value = 9
# Row 7
$ws.Range("G7").Value = 3
$ws.Range("H7").Value = 2.8
$ws.Range("I7").Value = 2.5
$ws.Range("J7").Value = 3.5
$ws.Range("L7").Value = 3.1
$ws.Range("N7").Value = 6.85
$ws.Range("O7").Value = 1.39
$ws.Range("P7").Value = 2.55
$ws.Range("S7").Value = 1.44
$ws.Range("T7").Value = 2.42
$ws.Range("U7").Value = 1.78
$ws.Range("V7").Value = 1.82
$ws.Range("Y7").Value = 10.5
$ws.Range("AA7").Value = 28
$ws.Range("AB7").Value = 37
$ws.Range("AC7").Value = 7.1
$ws.Range("AD7").Value = 5.5
$ws.Range("AG7").Value = 600
$ws.Range("AI7").Value = 11.75
$ws.Range("AN7").Value = 4.85
$ws.Range("AT7").Value = 2.4
$ws.Range("AU7").Value = 6.5
$ws.Range("AY7").Value = 21
# Row 8
$ws.Range("G8").Value = 1.8
$ws.Range("H8").Value = 3.5
$ws.Range("I8").Value = 4.33
$ws.Range("L8").Value = 4
$ws.Range("M8").Value = 1.05
$ws.Range("O8").Value = 1.25
$ws.Range("S8").Value = 1.33
$ws.Range("T8").Value = 3.25
$ws.Range("U8").Value = 1.67
$ws.Range("V8").Value = 2.1
$ws.Range("W8").Value = 8.5
$ws.Range("Y8").Value = 8.5
$ws.Range("Z8").Value = 15
$ws.Range("AA8").Value = 13
$ws.Range("AC8").Value = 12
$ws.Range("AD8").Value = 7
$ws.Range("AG8").Value = 151
$ws.Range("AH8").Value = 15
$ws.Range("AI8").Value = 23
$ws.Range("AJ8").Value = 15
$ws.Range("AL8").Value = 34
$ws.Range("AO8").Value = 9.5
$ws.Range("AP8").Value = 19
$ws.Range("AQ8").Value = 29
$ws.Range("AT8").Value = 3.25
$ws.Range("AU8").Value = 7.5
$ws.Range("AY8").Value = 26
$ws.Range("BB8").Value = 151
$ws.Range("BC8").Value = 151
$ws.Range("BD8").Value = 151
# Row 9
$ws.Range("J9").Value = 2.1
$ws.Range("K9").Value = 2.38
$ws.Range("L9").Value = 5
$ws.Range("M9").Value = 1.04
$ws.Range("O9").Value = 1.2
$ws.Range("Q9").Value = 1.7
$ws.Range("R9").Value = 2.1
$ws.Range("S9").Value = 1.29
$ws.Range("T9").Value = 3.5
$ws.Range("U9").Value = 1.67
$ws.Range("V9").Value = 2.1
$ws.Range("W9").Value = 9
$ws.Range("Z9").Value = 11
$ws.Range("AA9").Value = 11
$ws.Range("AB9").Value = 21
$ws.Range("AC9").Value = 17
$ws.Range("AD9").Value = 8.5
$ws.Range("AG9").Value = 151
$ws.Range("AH9").Value = 21
$ws.Range("AI9").Value = 34
$ws.Range("AJ9").Value = 19
$ws.Range("AK9").Value = 67
# Row 10
$ws.Range("G10").Value = 2.25
$ws.Range("H10").Value = 3.1
$ws.Range("J10").Value = 3
$ws.Range("K10").Value = 2
$ws.Range("M10").Value = 1.06
$ws.Range("O10").Value = 1.29
$ws.Range("Q10").Value = 2.35
$ws.Range("R10").Value = 1.57
$ws.Range("AW10").Value = 5
# Row 11
$ws.Range("R11").Value = 1.44
# Row 12
$ws.Range("G12").Value = 2.15
$ws.Range("H12").Value = 3.4
$ws.Range("K12").Value = 2
$ws.Range("L12").Value = 4.33
$ws.Range("Q12").Value = 2.35
$ws.Range("R12").Value = 1.57
$ws.Range("AG12").Value = 1000
$ws.Range("AN12").Value = 4
$ws.Range("AO12").Value = 13
$ws.Range("AQ12").Value = 41
$ws.Range("AX12").Value = 21
# Row 13
$ws.Range("G13").Value = 1.38
$ws.Range("M13").Value = 1.05
$ws.Range("N13").Value = 11
$ws.Range("O13").Value = 1.22
$ws.Range("P13").Value = 4
$ws.Range("R13").Value = 1.85
$ws.Range("S13").Value = 1.36
$ws.Range("T13").Value = 3
$ws.Range("W13").Value = 7
$ws.Range("Z13").Value = 8
$ws.Range("AA13").Value = 12
$ws.Range("AC13").Value = 11
$ws.Range("AD13").Value = 9.5
$ws.Range("AE13").Value = 23
$ws.Range("AH13").Value = 23
$ws.Range("AI13").Value = 51
$ws.Range("AJ13").Value = 29
$ws.Range("AK13").Value = 126
$ws.Range("AN13").Value = 3.25
$ws.Range("AT13").Value = 3
